{"js": "// Fix the ISA Design table:\n//  - sltR0 row: \"Otherwise, R0 -= 1\"  ->  \"Otherwise, R0 = 1\"\n//  - seqR0 row: \"Otherwise, R0 -= 1\"  ->  \"Otherwise, R0 = 1\"\n//  - beqR0 row: \"Otherwise, R0 -= 1\"  ->  \"Otherwise, PC++\"\n\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < rows.items.length; i++) {\n  rows.items[i].cells.load(\"items\");\n}\nawait context.sync();\n\n// Map instruction-name (first cell of each row) -> desired replacement text\n// for the last paragraph (\"R0 -= 1\") in the last cell (\"Function\" column).\nconst replacements = {\n  \"sltR0\": \"R0 = 1\",\n  \"seqR0\": \"R0 = 1\",\n  \"beqR0\": \"PC++\",\n};\n\n// Read the first cell of every row to identify which rows we need.\nfor (let i = 0; i < rows.items.length; i++) {\n  rows.items[i].cells.items[0].body.load(\"text\");\n}\nawait context.sync();\n\nconst targetRowIndexes = [];\nfor (let i = 0; i < rows.items.length; i++) {\n  const label = rows.items[i].cells.items[0].body.text.trim();\n  if (Object.prototype.hasOwnProperty.call(replacements, label)) {\n    targetRowIndexes.push(i);\n  }\n}\n\n// Load paragraphs of the Function (last) cell for each target row.\nconst lastCells = targetRowIndexes.map((idx) => {\n  const cells = rows.items[idx].cells.items;\n  return cells[cells.length - 1];\n});\nlastCells.forEach((cell) => cell.body.paragraphs.load(\"items/text\"));\nawait context.sync();\n\nfor (let k = 0; k < targetRowIndexes.length; k++) {\n  const idx = targetRowIndexes[k];\n  const label = rows.items[idx].cells.items[0].body.text.trim();\n  const newText = replacements[label];\n  const paras = lastCells[k].body.paragraphs.items;\n  for (let p = 0; p < paras.length; p++) {\n    if (paras[p].text.trim() === \"R0 -= 1\") {\n      paras[p].insertText(newText, \"Replace\");\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fix the ISA Design table (Word COM / PowerShell-style):\n#  - sltR0 row: \"Otherwise, R0 -= 1\"  ->  \"Otherwise, R0 = 1\"\n#  - seqR0 row: \"Otherwise, R0 -= 1\"  ->  \"Otherwise, R0 = 1\"\n#  - beqR0 row: \"Otherwise, R0 -= 1\"  ->  \"Otherwise, PC++\"\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\n# Instruction name (first column) -> new text for the \"Function\" (last) column\n$replacements = @{\n    \"sltR0\" = \"R0 = 1\"\n    \"seqR0\" = \"R0 = 1\"\n    \"beqR0\" = \"PC++\"\n}\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    # Cell text always carries a trailing cell-mark (CR + BEL); strip it off.\n    $label = $table.Cell($r, 1).Range.Text.TrimEnd([char]13, [char]7)\n    if ($replacements.ContainsKey($label)) {\n        $newText = $replacements[$label]\n        $cellRange = $table.Cell($r, $colCount).Range\n        # wdFindStop (0) keeps the search confined to this range;\n        # wdReplaceOne (1) replaces just the single match found.\n        $cellRange.Find.ClearFormatting()\n        $cellRange.Find.Execute(\"R0 -= 1\", $false, $false, $false, $false, $false, $true, 0, $false, $newText, 1) | Out-Null\n    }\n}\n"}
